# Auto-generated Excel COM-interop script
# Applies market-price / profit-column refresh to the Leve Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1959.6
$ws.Range("I2").Value = 800
$ws.Range("J2").Value = 2249.5
$ws.Range("K2").Value = 800
$ws.Range("L2").Value = 2249.5
$ws.Range("M2").Value = -687
$ws.Range("N2").Value = -2475.5

$ws.Range("H4").Value = 1151.6818
$ws.Range("I4").Value = 723.5454999999999
$ws.Range("J4").Value = 1579.8182
$ws.Range("K4").Value = 723.5454999999999
$ws.Range("L4").Value = 1579.8182
$ws.Range("M4").Value = -609.5454999999999
$ws.Range("N4").Value = -1807.8182

$ws.Range("H8").Value = 2389.25
$ws.Range("I8").Value = 2389.25
$ws.Range("K8").Value = 7167.75
$ws.Range("M8").Value = -7028.75

$ws.Range("H9").Value = 293.55554
$ws.Range("I9").Value = 252.33333
$ws.Range("K9").Value = 252.33333
$ws.Range("M9").Value = -83.33332999999999

$ws.Range("H12").Value = 660.2222
$ws.Range("I12").Value = 423.14285
$ws.Range("J12").Value = 811.0909
$ws.Range("K12").Value = 423.14285
$ws.Range("L12").Value = 811.0909
$ws.Range("M12").Value = -253.14285
$ws.Range("N12").Value = -1151.0909

$ws.Range("H33").Value = 548.0833
$ws.Range("J33").Value = 683
$ws.Range("L33").Value = 683
$ws.Range("N33").Value = -1141

$ws.Range("H55").Value = 305.16666
$ws.Range("J55").Value = 398.7143
$ws.Range("L55").Value = 398.7143
$ws.Range("N55").Value = -826.7143

$ws.Range("H80").Value = 2588756.5
$ws.Range("I80").Value = 1636435.6
$ws.Range("J80").Value = 3405031.2
$ws.Range("K80").Value = 4909306.800000001
$ws.Range("L80").Value = 10215093.6
$ws.Range("M80").Value = -4908308.800000001
$ws.Range("N80").Value = -10217089.6

$ws.Range("H83").Value = 2588756.5
$ws.Range("I83").Value = 1636435.6
$ws.Range("J83").Value = 3405031.2
$ws.Range("K83").Value = 14727920.4
$ws.Range("L83").Value = 30645280.8
$ws.Range("M83").Value = -14722928.4
$ws.Range("N83").Value = -30655264.8

$ws.Range("H132").Value = 2293.311
$ws.Range("I132").Value = 2324.4524
$ws.Range("J132").Value = 1857.3334
$ws.Range("K132").Value = 6973.3572
$ws.Range("L132").Value = 5572.0002
$ws.Range("M132").Value = -4443.3572
$ws.Range("N132").Value = -10632.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 652.617
$ws.Range("I2").Value = 564.94446
$ws.Range("K2").Value = 564.94446
$ws.Range("M2").Value = -451.94446

$ws.Range("H19").Value = 2282.6667
$ws.Range("I19").Value = 2449
$ws.Range("J19").Value = 1950
$ws.Range("K19").Value = 2449
$ws.Range("L19").Value = 1950
$ws.Range("M19").Value = -2220
$ws.Range("N19").Value = -2408

$ws.Range("H30").Value = 333370000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 333370000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 333370000
$ws.Range("N30").Value = -333370300
$ws.Range("M30").ClearContents()

$ws.Range("H32").Value = 13302.425
$ws.Range("I32").Value = 13592.079
$ws.Range("K32").Value = 13592.079
$ws.Range("M32").Value = -13305.079

$ws.Range("H74").Value = 2873.5
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 2873.5
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H97").Value = 2852.0952
$ws.Range("I97").Value = 1665.5555
$ws.Range("K97").Value = 1665.5555
$ws.Range("M97").Value = -1169.5555

$ws.Range("H102").Value = 2447
$ws.Range("I102").Value = 1590.6316
$ws.Range("K102").Value = 1590.6316
$ws.Range("M102").Value = 31.36840000000007

$ws.Range("H110").Value = 8150
$ws.Range("I110").Value = 10437.875
$ws.Range("J110").Value = 5535.2856
$ws.Range("K110").Value = 10437.875
$ws.Range("L110").Value = 5535.2856
$ws.Range("M110").Value = -8392.875
$ws.Range("N110").Value = -9625.285599999999

$ws.Range("H116").Value = 652.617
$ws.Range("I116").Value = 564.94446
$ws.Range("K116").Value = 564.94446
$ws.Range("M116").Value = 1729.05554

$ws.Range("H132").Value = 7698021
$ws.Range("I132").Value = 6476
$ws.Range("J132").Value = 33336504
$ws.Range("K132").Value = 19428
$ws.Range("L132").Value = 100009512
$ws.Range("M132").Value = -16898
$ws.Range("N132").Value = -100014572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 652.617
$ws.Range("I3").Value = 564.94446
$ws.Range("K3").Value = 564.94446
$ws.Range("M3").Value = -450.94446

$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws.Range("H86").Value = 6099
$ws.Range("I86").Value = 5148.75
$ws.Range("K86").Value = 5148.75
$ws.Range("M86").Value = -4025.75

$ws.Range("H89").Value = 6099
$ws.Range("I89").Value = 5148.75
$ws.Range("K89").Value = 25743.75
$ws.Range("M89").Value = -20127.75

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H134").Value = 11112360
$ws.Range("J134").Value = 25001350
$ws.Range("L134").Value = 75004050
$ws.Range("N134").Value = -75009120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 11112674
$ws.Range("I16").Value = 14286581
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 14286581
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = -14286294
$ws.Range("N16").Value = -4574

$ws.Range("H58").Value = 5742.7
$ws.Range("J58").Value = 8602.6
$ws.Range("L58").Value = 8602.6
$ws.Range("N58").Value = -9008.6

$ws.Range("H103").Value = 32829
$ws.Range("I103").Value = 12862.667
$ws.Range("K103").Value = 12862.667
$ws.Range("M103").Value = -11690.667

$ws.Range("H113").Value = 11112674
$ws.Range("I113").Value = 14286581
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 14286581
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -14284411
$ws.Range("N113").Value = -8340

$ws.Range("H136").Value = 5742.7
$ws.Range("J136").Value = 8602.6
$ws.Range("L136").Value = 25807.8
$ws.Range("N136").Value = -30907.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 7787.3184
$ws.Range("I134").Value = 2666
$ws.Range("K134").Value = 7998
$ws.Range("M134").Value = -2928

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 9000
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H21").Value = 50021250
$ws.Range("I21").Value = 66685664
$ws.Range("K21").Value = 66685664
$ws.Range("M21").Value = -66685491

$ws.Range("H30").Value = 50021250
$ws.Range("I30").Value = 66685664
$ws.Range("K30").Value = 66685664
$ws.Range("M30").Value = -66685559

$ws.Range("H126").Value = 18116900
$ws.Range("I126").Value = 25081674
$ws.Range("J126").Value = 8484.799999999999
$ws.Range("K126").Value = 75245022
$ws.Range("L126").Value = 25454.4
$ws.Range("M126").Value = -75242552
$ws.Range("N126").Value = -30394.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 757.2727
$ws.Range("I22").Value = 687.3333
$ws.Range("J22").Value = 841.2
$ws.Range("K22").Value = 687.3333
$ws.Range("L22").Value = 841.2
$ws.Range("M22").Value = -392.3333
$ws.Range("N22").Value = -1431.2

$ws.Range("H27").Value = 757.2727
$ws.Range("I27").Value = 687.3333
$ws.Range("J27").Value = 841.2
$ws.Range("K27").Value = 687.3333
$ws.Range("L27").Value = 841.2
$ws.Range("M27").Value = -580.3333
$ws.Range("N27").Value = -1055.2

$ws.Range("H46").Value = 864.3077
$ws.Range("I46").Value = 693.2222
$ws.Range("K46").Value = 693.2222
$ws.Range("M46").Value = -505.2222

$ws.Range("H55").Value = 998.9091
$ws.Range("I55").Value = 431.8
$ws.Range("K55").Value = 431.8
$ws.Range("M55").Value = -258.8

$ws.Range("H68").Value = 2454420
$ws.Range("I68").Value = 4632438
$ws.Range("J68").Value = 4149.75
$ws.Range("K68").Value = 4632438
$ws.Range("L68").Value = 4149.75
$ws.Range("M68").Value = -4631689
$ws.Range("N68").Value = -5647.75

$ws.Range("H71").Value = 2454420
$ws.Range("I71").Value = 4632438
$ws.Range("J71").Value = 4149.75
$ws.Range("K71").Value = 23162190
$ws.Range("L71").Value = 20748.75
$ws.Range("M71").Value = -23158446
$ws.Range("N71").Value = -28236.75

$ws.Range("H136").Value = 4053.4827
$ws.Range("I136").Value = 2694.0667
$ws.Range("K136").Value = 8082.2001
$ws.Range("M136").Value = -5532.2001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 60249.832
$ws.Range("I2").Value = 67999.75
$ws.Range("J2").Value = 44750
$ws.Range("K2").Value = 67999.75
$ws.Range("L2").Value = 44750
$ws.Range("M2").Value = -67887.75
$ws.Range("N2").Value = -44974

$ws.Range("H81").Value = 1759.3636
$ws.Range("I81").Value = 1899.125
$ws.Range("J81").Value = 1386.6666
$ws.Range("K81").Value = 3798.25
$ws.Range("L81").Value = 2773.3332
$ws.Range("M81").Value = -2737.25
$ws.Range("N81").Value = -4895.3332

$ws.Range("H84").Value = 1759.3636
$ws.Range("I84").Value = 1899.125
$ws.Range("J84").Value = 1386.6666
$ws.Range("K84").Value = 18991.25
$ws.Range("L84").Value = 13866.666
$ws.Range("M84").Value = -13687.25
$ws.Range("N84").Value = -24474.666

$ws.Range("H136").Value = 281911.4
$ws.Range("I136").Value = 4399.517
$ws.Range("J136").Value = 1431603.6
$ws.Range("K136").Value = 13198.551
$ws.Range("L136").Value = 4294810.800000001
$ws.Range("M136").Value = -10648.551
$ws.Range("N136").Value = -4299910.800000001
